$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.358.17"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "'1.941.51"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'251.93"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'0.7223"
$ws.Range("E6").Value = "  -8.35%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.3345"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").Value = "'28.83"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").Value = "'0.07414"
$ws.Range("E10").Value = "  +6.05%  "
$ws.Range("D11").Value = "'0.8178"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "'0.08141"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "'1.941.46"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "'5.493"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "'95.29"
$ws.Range("E15").Value = "  -4.89%  "
$ws.Range("D16").Value = "'14.94"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000008451"
$ws.Range("E17").Value = "  +6.73%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'30.380.61"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'253.86"
$ws.Range("E19").Value = "  -7.01%  "
$ws.Range("D20").Value = "'5.897"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "'2.198.82"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'6.990"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'9.866"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "'163.24"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").Value = "'2.417"
$ws.Range("E27").Value = "  +4.50%  "
$ws.Range("D28").Value = "'19.37"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").Value = "'0.1322"
$ws.Range("E29").Value = "  -11.38%  "
$ws.Range("D30").Value = "'1.576"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "'4.468"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").Value = "'4.257"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("D34").Value = "'0.05295"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "'1.308"
$ws.Range("E35").Value = "  +7.04%  "
$ws.Range("D36").Value = "'0.7582"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "'2.745"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").Value = "'0.01993"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'2.852"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "'81.19"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").Value = "'6.618"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.4577"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").Value = "'0.8474"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'102.95"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "'9.818"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'7.523"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "'36.91"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "'0.4203"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'1.512"
$ws.Range("E51").Value = "  -1.17%  "
